$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 185, shifting existing rows 185:265 down to 186:266
$ws.Rows("185:185").Insert()

# Populate the newly inserted row 185 with the new data record
$ws.Range("A185").Value = 4
$ws.Range("B185").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C185").Value = "Los Lagos"
$ws.Range("D185").Value = 44704
$ws.Range("E185").Value = 10
$ws.Range("F185").Value = 100112003
$ws.Range("G185").Value = "Ajo"
$ws.Range("H185").Value = "Chino"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 60
$ws.Range("K185").Value = 22000
$ws.Range("L185").Value = 22000
$ws.Range("M185").Value = 22000
$ws.Range("N185").Value = "$/caja 10 kilos"
$ws.Range("O185").Value = "China"
$ws.Range("P185").Value = 2200
$ws.Range("Q185").Value = 10
$ws.Range("R185").Value = "Hortaliza"
